# Update the "Förändrad" (Changed) date column (C) for all data rows
# from 2023-10-05 (serial 45204) to 2023-10-06 (serial 45205).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlUp = -4162
$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End($xlUp).Row
if ($lastRow -lt 2) { $lastRow = 89 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45204) {
        $cell.Value = 45205
    }
}
